# Slide 2, shape 4 ("Tijdelijke aanduiding voor tekst 6" / ph type="body" idx="14")
# currently just contains the date "27 april 2021".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(4)

# Give the placeholder an explicit position/size (previously inherited from the
# layout, i.e. empty <p:spPr/>). Values are in points; PowerPoint COM stores
# Left/Top/Width/Height in points (1 pt = 12700 EMU).
$shp.Left = 113.43331146240234
$shp.Top = 353.8030090332031
$shp.Width = 429.5391540527344
$shp.Height = 127.37504577636719

# Append a new paragraph with the email address after the existing date line.
$tr = $shp.TextFrame.TextRange
[void]$tr.InsertAfter("`r" + "dmitriy.vanderelst@student.arteveldehs.be")
